$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.099.92'
$ws.Range('E2').Value = '  +1.88%  '
$ws.Range('D3').Value = '3.826.33'
$ws.Range('E3').Value = '  +0.43%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.993'
$ws.Range('E4').Value = '  -0.34%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '626.81'
$ws.Range('E5').Value = '  +4.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '165.24'
$ws.Range('D7').Value = '3.823.33'
$ws.Range('E7').Value = '  +0.42%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  +0.28%  '
$ws.Range('E10').Value = '  +0.97%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.455'
$ws.Range('E11').Value = '  +0.31%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.66'
$ws.Range('E12').Value = '  +4.20%  '
$ws.Range('E13').Value = '  +0.06%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.11'
$ws.Range('E14').Value = '  +0.11%  '
$ws.Range('D15').Value = '4.468.95'
$ws.Range('E15').Value = '  +0.36%  '
$ws.Range('D16').Value = '3.860.83'
$ws.Range('E16').Value = '  +1.11%  '
$ws.Range('D17').Value = '69.091.44'
$ws.Range('E17').Value = '  +1.81%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.34'
$ws.Range('E18').Value = '  -0.97%  '
$ws.Range('E19').Value = '  +0.98%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.114'
$ws.Range('E20').Value = '  +0.09%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '468.13'
$ws.Range('E21').Value = '  +0.87%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.72'
$ws.Range('E22').Value = '  -1.31%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.710'
$ws.Range('E23').Value = '  +0.83%  '
$ws.Range('E24').Value = '  +3.64%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.14'
$ws.Range('E25').Value = '  +1.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.08'
$ws.Range('E26').Value = '  -0.36%  '
$ws.Range('E27').Value = '  +1.68%  '
$ws.Range('E28').Value = '  +0.08%  '
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('D30').Value = '3.980.01'
$ws.Range('E30').Value = '  +0.51%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.26'
$ws.Range('E31').Value = '  +2.09%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.67'
$ws.Range('E32').Value = '  -4.06%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.36'
$ws.Range('E33').Value = '  -1.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '29.22'
$ws.Range('E34').Value = '  -0.75%  '
$ws.Range('E35').Value = '  +0.57%  '
$ws.Range('E36').Value = '  -0.19%  '
$ws.Range('E37').Value = '  +1.79%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.149'
$ws.Range('E38').Value = '  +7.55%  '
$ws.Range('B39').Value = 'dogwifhat'
$ws.Range('C39').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.32'
$ws.Range('E39').Value = '  +2.14%  '
$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.93'
$ws.Range('E40').Value = '  +2.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.982'
$ws.Range('E41').Value = '  -1.29%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '156.10'
$ws.Range('E44').Value = '  +3.22%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.302'
$ws.Range('E45').Value = '  +0.36%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.42'
$ws.Range('E46').Value = '  +1.79%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '46.79'
$ws.Range('E47').Value = '  -2.14%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '42.89'
$ws.Range('E48').Value = '  -5.85%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.46'
$ws.Range('E49').Value = '  +0.97%  '
$ws.Range('E50').Value = '  +2.01%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '382.12'
$ws.Range('E51').Value = '  -2.96%  '
